$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.960.94"
$ws.Range("E2").Value = "  +4.75%  "

$ws.Range("D3").Value = "3.084.93"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'581.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "

$ws.Range("D6").Value = "'142.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.39%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.076.02"
$ws.Range("E8").Value = "  +3.26%  "

$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("E10").Value = "  +5.15%  "

$ws.Range("D11").Value = "'5.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.79%  "

$ws.Range("E12").Value = "  +2.80%  "

$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("D14").Value = "'35.43"
$ws.Range("D14").Style = "Normal"

$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").Value = "3.592.05"
$ws.Range("E16").Value = "  +3.08%  "

$ws.Range("E17").Value = "  +3.04%  "

$ws.Range("D18").Value = "3.079.14"
$ws.Range("E18").Value = "  +3.00%  "

$ws.Range("D19").Value = "61.888.32"
$ws.Range("E19").Value = "  +4.72%  "

$ws.Range("D20").Value = "'448.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.68%  "

$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("D22").Value = "'0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.61%  "

$ws.Range("E23").Value = "  +5.05%  "

$ws.Range("D24").Value = "'13.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.86%  "

$ws.Range("D25").Value = "'82.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  +6.47%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  +5.09%  "

$ws.Range("D30").Value = "'8.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.54%  "

$ws.Range("D31").Value = "'6.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.16%  "

$ws.Range("D32").Value = "'0.113"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.71%  "

$ws.Range("D33").Value = "'26.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.48%  "

$ws.Range("E34").Value = "  +4.71%  "

$ws.Range("E35").Value = "  +2.59%  "

$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("E37").Value = "  +5.75%  "

$ws.Range("D38").Value = "'50.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("E39").Value = "  +10.27%  "

$ws.Range("E40").Value = "  +2.37%  "

$ws.Range("D41").Value = "'426.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.92%  "

$ws.Range("E42").Value = "  +5.84%  "

$ws.Range("D43").Value = "2.910.96"
$ws.Range("E43").Value = "  +4.78%  "

$ws.Range("E44").Value = "  +8.16%  "

$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").Value = "'2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.89%  "

$ws.Range("E47").Value = "  +4.87%  "

$ws.Range("D49").Value = "'124.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").Value = "'24.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.04%  "
